$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "40.036.47"
$ws.Range("E2").Value = "  +0.91%  "

# Row 3
$ws.Range("D3").Value = "2.216.93"
$ws.Range("E3").Value = "  -0.02%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").Value = "'290.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.76%  "

# Row 6
$ws.Range("D6").Value = "'88.48"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.64%  "

# Row 7
$ws.Range("E7").Value = "  +0.05%  "

# Row 8
$ws.Range("E8").Value = "  -0.02%  "

# Row 9
$ws.Range("E9").Value = "  +1.23%  "

# Row 10
$ws.Range("D10").Value = "'30.81"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.87%  "

# Row 11
$ws.Range("D11").Value = "'0.0782"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.57%  "

# Row 12
$ws.Range("D12").Value = "'47.92"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.78%  "

# Row 13
$ws.Range("E13").Value = "  +2.61%  "

# Row 14
$ws.Range("D14").Value = "'6.47"
$ws.Range("D14").Style = "Normal"

# Row 15
$ws.Range("D15").Value = "2.560.65"
$ws.Range("E15").Value = "  +0.10%  "

# Row 16
$ws.Range("D16").Value = "'14.02"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.79%  "

# Row 17
$ws.Range("D17").Value = "2.209.72"
$ws.Range("E17").Value = "  -0.14%  "

# Row 18
$ws.Range("E18").Value = "  +1.56%  "

# Row 19
$ws.Range("D19").Value = "39.985.84"
$ws.Range("E19").Value = "  +0.99%  "

# Row 20
$ws.Range("D20").Value = "'11.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +13.98%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0885"
$ws.Range("E21").Value = "  +0.71%  "

# Row 22
$ws.Range("D22").Value = "'5.82"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.38%  "

# Row 23
$ws.Range("D23").Value = "'65.62"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.91%  "

# Row 24
$ws.Range("D24").Value = "'235.63"
$ws.Range("D24").Style = "Normal"

# Row 25
$ws.Range("E25").Value = "  +0.07%  "

# Row 26
$ws.Range("D26").Value = "'2.46"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.24%  "

# Row 27
$ws.Range("E27").Value = "  +0.11%  "

# Row 28
$ws.Range("D28").Value = "'22.63"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.50%  "

# Row 29
$ws.Range("E29").Value = "  +0.96%  "

# Row 30
$ws.Range("D30").Value = "'9.24"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.80%  "

# Row 31
$ws.Range("D31").Value = "'153.07"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.66%  "

# Row 32
$ws.Range("D32").Value = "'32.27"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.07%  "

# Row 34
$ws.Range("D34").Value = "'4.98"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.78%  "

# Row 35
$ws.Range("D35").Value = "'0.0721"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.63%  "

# Row 36
$ws.Range("E36").Value = "  +0.13%  "

# Row 37
$ws.Range("E37").Value = "  +7.10%  "

# Row 38
$ws.Range("E38").Value = "  +0.56%  "

# Row 39
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "'0.100"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.06%  "

# Row 40
$ws.Range("B40").Value = "Celestia"
$ws.Range("C40").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D40").Value = "'15.89"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.46%  "

# Row 41
$ws.Range("E41").Value = "  +3.32%  "

# Row 42
$ws.Range("D42").Value = "2.100.98"

# Row 43
$ws.Range("E43").Value = "  +4.71%  "

# Row 44
$ws.Range("E44").Value = "  +1.85%  "

# Row 45
$ws.Range("E45").Value = "  +1.08%  "

# Row 46
$ws.Range("D46").Value = "'17.71"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +9.15%  "

# Row 47
$ws.Range("E47").Value = "  +7.06%  "

# Row 48
$ws.Range("E48").Value = "  +1.71%  "

# Row 49
$ws.Range("D49").Value = "2.432.22"
$ws.Range("E49").Value = "  +0.11%  "

# Row 50
$ws.Range("D50").Value = "'69.54"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.89%  "

# Row 51
$ws.Range("D51").Value = "'88.75"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.03%  "
